$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: update C2 text and move the saved selection
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C2").Value = "Cash. 046C000350 - Cá nhân trong nu?c 2 Cá nhân trong nu?c 2 Cá nhân trong nu?c 2"
$ws1.Range("C11").Select()

# ---------------------------------------------------------------------------
# Sheet3: insert two new leading columns (username/pass) and a trailing
# column (tradingPass), fill in the new login fields, tweak a couple of the
# existing order values, and restyle the new header cells to match the
# other header columns (fill only, no bold / mono font).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

# Shift the existing af/orderType/symbol/qtty/price/BS/priceType columns
# from A:G to C:I, preserving their values/styles/widths.
$ws3.Range("A1:B1").EntireColumn.Insert()

# New header cells - copy the fill-only header style (used elsewhere in the
# workbook) from Sheet1!A1, then strip the bold so it matches fontId 0.
$ws1.Range("A1").Copy()
$ws3.Range("A1").PasteSpecial(-4122)
$ws3.Range("B1").PasteSpecial(-4122)
$ws3.Range("J1").PasteSpecial(-4122)
$ws3.Range("A1,B1,J1").Font.Bold = $false

$ws3.Range("A1").Value = "username"
$ws3.Range("B1").Value = "pass"
$ws3.Range("J1").Value = "tradingPass"

# New data row (row 2) - login fields.
$ws3.Range("A2").Value = "046FIA0016"
$ws3.Range("B2").Value = 123
# Leading-zero account number must stay text; the leading apostrophe mimics
# typing it in Excel, which also picks up the existing quote-prefix style.
$ws3.Range("C2").Value = "'0001000021"
$ws3.Range("J2").Value = 123

# Updated order values.
$ws3.Range("G2").Value = 25
$ws3.Range("H2").Value = "Mua"

# Column widths (A:J).
$ws3.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws3.Columns.Item(2).ColumnWidth = 14.666666666666666
$ws3.Columns.Item(3).ColumnWidth = 12.666666666666666
$ws3.Columns.Item(4).ColumnWidth = 20
$ws3.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws3.Columns.Item(6).ColumnWidth = 12.166666666666666
$ws3.Columns.Item(7).ColumnWidth = 11
$ws3.Columns.Item(8).ColumnWidth = 9.666666666666666
$ws3.Columns.Item(9).ColumnWidth = 14.5
$ws3.Columns.Item(10).ColumnWidth = 12.166666666666666

$ws3.Range("C6").Select()
